$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date in column C for rows 2-15 to 2023-09-14 (serial 45183)
$newDate = Get-Date -Year 2023 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0
foreach ($row in 2..15) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
